# Updated file with local changes
#
# The original sheet had:
#   A1 = 0              (bold, bordered, centered style)
#   A2 = long JSON-ish "questions = [...]" string (plain style)
#
# The edit:
#   - removes row 2 entirely (its text now lives in A1)
#   - puts the (re-formatted / pretty-printed) questions text into A1
#   - drops the now-unused bold/border formatting that used to live on A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Get rid of the old row 2 (the row that used to hold the questions text).
# Deleting the whole row shifts everything up and shrinks the used range
# back down to just row 1.
$ws.Rows.Item(2).Delete()

# A1 used to be a bold, centered, bordered "0" placeholder cell -- clear
# that formatting so the cell goes back to the workbook's default style.
$ws.Range("A1").ClearFormats()

$questionsText = @'
questions = [
    {
        "title": "You are a data engineer for a global company that collects data from various sources, such as IoT devices, customer interactions, and third-party providers. The data comes in multiple formats like CSV, JSON, and XML. You want to ingest this data efficiently into the AWS environment for later analysis.Which AWS services should you use?",
        "ques_type": 2,
        "options": [
            "Use Amazon Kinesis Data Streams to ingest data in real-time and Amazon Simple Storage Service (S3) for storing different data formats.",
            "Use AWS Data Pipeline to ingest data and AWS DynamoDB for storage.",
            "Use Amazon Simple Queue Service (SQS) for data ingestion and Amazon Relational Database Service (RDS) for storage.",
            "Use Amazon Database Migration Service (DMS) to import data and Amazon Elastic Block Store (EBS) for storage."
        ],
        "score": "Use Amazon Kinesis Data Streams to ingest data in real-time and Amazon Simple Storage Service (S3) for storing different data formats."
    },
    {
        "title": "You are a data analyst and have been tasked with cleaning and enriching a large dataset collected from various e-commerce sites. The dataset contains missing values, duplicates, and outliers. It also requires enrichment with demographic data for a more comprehensive analysis.Which AWS services should you use?",
        "ques_type": 2,
        "options": [
            "Use AWS Glue for data cleaning and Amazon Simple Storage Service (S3) Select for data enrichment.",
            "Use AWS Database Migration Service (DMS) for data cleaning and Amazon Redshift for data enrichment.",
            "Use Amazon Athena for data cleaning and AWS Lambda for data enrichment.",
            "Use Amazon Relational Database Service (RDS) for data cleaning and AWS Glue for data enrichment."
        ],
        "score": "Use AWS Glue for data cleaning and Amazon Simple Storage Service (S3) Select for data enrichment."
    },
    {
        "title": "You are working as a data analyst in a company that has a significant amount of data stored in Amazon Simple Storage Service (S3). The data is in CSV format and is updated on an hourly basis. You want to analyze this data to derive insights.What should you do?",
        "ques_type": 2,
        "options": [
            "Load the CSV data into a Pandas DataFrame and perform the analysis using Python.",
            "Run SQL-like queries directly on the data using a suitable AWS service.",
            "Import the data into an Amazon Redshift cluster and run the analysis.",
            "Use AWS Glue to transform the CSV data into a more suitable format for analysis."
        ],
        "score": "Run SQL-like queries directly on the data using a suitable AWS service."
    },
    {
        "title": "You are an AWS Solutions Architect at a large tech company. Your system uses Amazon Relational Database Service (RDS) for PostgreSQL. Some specific SELECT queries are taking a significant amount of time to return results, affecting overall system performance.What should you do?",
        "ques_type": 2,
        "options": [
            "Create multi-column indexes covering the queried fields.",
            "Increase the instance type of your RDS instance.",
            "Use Amazon Redshift for query processing.",
            "Use a Read Replica for SELECT operations."
        ],
        "score": "Create multi-column indexes covering the queried fields."
    }
]
'@

# The here-string captures a trailing newline -- trim it so the cell value
# matches the original text exactly (no trailing blank line).
$questionsText = $questionsText.TrimEnd("`r", "`n")

$ws.Range("A1").Value = $questionsText

# Writing a multi-line value auto-expands the row height with an explicit
# customHeight; AutoFit() puts the row back to the sheet's implicit default
# height (matching the original, which had no per-row height override).
$ws.Rows.Item(1).AutoFit()
